# Fix format of SYNoEVC file
#
# The "Start Year Number of EV Chargers" (SYNoEVC) sheet used to show the
# full trajectory from 2021-2050 across row 1 (years) / row 2 (values).
# It is being reformatted to a single "start year" input: only column B
# is kept (now showing 2020 instead of 2021, and pulling its value from
# Calculations!B8 instead of Calculations!C11); every other year/value
# cell in columns C:AE is cleared out (contents only - formatting is left
# untouched so the alternating column shading survives).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYNoEVC")

# Row 1 (years): B1 becomes the new start year, rest of the row is cleared.
$ws.Range("B1").Value = 2020
$ws.Range("C1:AE1").ClearContents()

# Row 2 (values): B2's formula now references the updated Calculations
# layout; rest of the row is cleared.
$ws.Range("B2").Formula = "=Calculations!B8"
$ws.Range("C2:AE2").ClearContents()

# The workbook was left with SYNoEVC as the active/selected sheet, with
# B3 selected (instead of the "About" sheet that used to be active).
$ws.Activate()
$ws.Range("B3").Select()
